{"js": "// Word Software Requirements Specification \u2014 update \"Druhy u\u017eivatel\u016f\",\n// \"Pou\u017eit\u00e9 technologie\" and the admin/moderator bullet list paragraphs.\n\nfunction wrapOoxml(bodyXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n</Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n${bodyXml}\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\n// 1) \"Aplikace rozli\u0161uje dva druhy u\u017eivatel\u016f. ...\" -> new 4-role description.\nconst oldUserTypes =\n  \"Aplikace rozli\u0161uje dva druhy u\u017eivatel\u016f. B\u011b\u017en\u00e9 n\u00e1v\u0161t\u011bvn\u00edky bez u\u017eivatelsk\u00e9ho \u00fa\u010dtu, kte\u0159\u00ed mohou aplikace b\u011b\u017en\u011b vyu\u017e\u00edvat. Administr\u00e1tory, kte\u0159\u00ed mohou spravovat po\u017eadavky na p\u0159id\u00e1n\u00ed postav od b\u011b\u017en\u00fdch u\u017eivatel\u016f, upravovat dosavadn\u00ed z\u00e1znamy o postav\u00e1ch a mazat je.\";\nconst newUserTypes =\n  \"Aplikace rozli\u0161uje 4 druhy u\u017eivatel\u016f. B\u011b\u017en\u00fd u\u017eivatel (p\u0159ihl\u00e1\u0161en\u00fd \u010di nep\u0159ihl\u00e1\u0161en\u00fd) m\u00e1 k dispozici pouze z\u00e1kladn\u00ed funkce. Ov\u011b\u0159en\u00fd u\u017eivatel m\u016f\u017ee nav\u00edc schvalovat po\u017eadavky na p\u0159id\u00e1n\u00ed postav. Moder\u00e1tor m\u016f\u017ee nav\u00edc postavy upravovat a mazat. Administr\u00e1tor m\u016f\u017ee nav\u00edc spravovat u\u017eivatele (mazat a m\u011bnit jejich \u00farovn\u011b).\";\n\nlet results = context.document.body.search(oldUserTypes, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(newUserTypes, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Aplikace je tvo\u0159ena v HTML a CSS, ...\" -> split into two paragraphs\n//    describing the Vue.js frontend / Java backend / Pinia auth.\nconst oldTech =\n  \"Aplikace je tvo\u0159ena v HTML a CSS, na frontend vyu\u017e\u00edv\u00e1 jazyk JavaScript a backend je v jazyku Java. Aplikace tak\u00e9 vyu\u017e\u00edv\u00e1 vlastn\u00ed MySQL datab\u00e1ze. Pro CSS je pou\u017eit framework Bootstrap.\";\n\nresults = context.document.body.search(oldTech, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  const techBodyXml =\n    '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Frontend</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> aplikace je tvo\u0159en pomoc\u00ed vue.js s pou\u017eit\u00edm frameworku </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>bootstrap</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">. </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Backend</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> aplikace je v jazyce Java. Aplikace tak\u00e9 vyu\u017e\u00edv\u00e1 vlastn\u00ed </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>MySQL</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> datab\u00e1ze.</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Na spr\u00e1vu p\u0159ihl\u00e1\u0161en\u00ed je pou\u017eit </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Pinia</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> s JWT tokeny.</w:t></w:r></w:p>';\n  results.items[0].insertOoxml(wrapOoxml(techBodyXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"Administr\u00e1to\u0159i mohou schvalovat/zam\u00edtnout ...\" -> now also ov\u011b\u0159en\u00ed\n//    u\u017eivatel\u00e9 a moder\u00e1to\u0159i.\nconst oldApprove =\n  \"Administr\u00e1to\u0159i mohou schvalovat/zam\u00edtnout jednotliv\u00e9 po\u017eadavky na p\u0159id\u00e1n\u00ed postav od u\u017eivatel\u016f.\";\nconst newApprove =\n  \"Administr\u00e1to\u0159i, ov\u011b\u0159en\u00ed u\u017eivatel\u00e9 a moder\u00e1to\u0159i mohou schvalovat/zam\u00edtnout jednotliv\u00e9 po\u017eadavky na p\u0159id\u00e1n\u00ed postav od u\u017eivatel\u016f.\";\n\nresults = context.document.body.search(oldApprove, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(newApprove, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \"Administr\u00e1to\u0159i mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je\n//    mazat.\" -> now also moder\u00e1to\u0159i, plus a brand-new following paragraph\n//    about administrators managing users.\nconst oldEdit =\n  \"Administr\u00e1to\u0159i mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je mazat.\";\n\nresults = context.document.body.search(oldEdit, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  const editBodyXml =\n    '<w:p><w:r><w:t xml:space=\"preserve\">Administr\u00e1to\u0159i </w:t></w:r><w:r><w:t xml:space=\"preserve\">a moder\u00e1to\u0159i </w:t></w:r><w:r><w:t>mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je mazat.</w:t></w:r></w:p>' +\n    \"<w:p><w:r><w:t>Administr\u00e1to\u0159i mohou upravit a mazat jednotliv\u00e9 u\u017eivatele.</w:t></w:r></w:p>\";\n  results.items[0].insertOoxml(wrapOoxml(editBodyXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word Software Requirements Specification \u2014 update \"Druhy u\u017eivatel\u016f\",\n# \"Pou\u017eit\u00e9 technologie\" and the admin/moderator bullet list paragraphs.\n\n$d = $word.ActiveDocument\n\n# 1) \"Aplikace rozli\u0161uje dva druhy u\u017eivatel\u016f. ...\" -> new 4-role description.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Aplikace rozli\u0161uje dva druhy u\u017eivatel\u016f. B\u011b\u017en\u00e9 n\u00e1v\u0161t\u011bvn\u00edky bez u\u017eivatelsk\u00e9ho \u00fa\u010dtu, kte\u0159\u00ed mohou aplikace b\u011b\u017en\u011b vyu\u017e\u00edvat. Administr\u00e1tory, kte\u0159\u00ed mohou spravovat po\u017eadavky na p\u0159id\u00e1n\u00ed postav od b\u011b\u017en\u00fdch u\u017eivatel\u016f, upravovat dosavadn\u00ed z\u00e1znamy o postav\u00e1ch a mazat je.\"\n$found = $find.Execute()\nif ($found) {\n    $rng = $find.Parent\n    $rng.Text = \"Aplikace rozli\u0161uje 4 druhy u\u017eivatel\u016f. B\u011b\u017en\u00fd u\u017eivatel (p\u0159ihl\u00e1\u0161en\u00fd \u010di nep\u0159ihl\u00e1\u0161en\u00fd) m\u00e1 k dispozici pouze z\u00e1kladn\u00ed funkce. Ov\u011b\u0159en\u00fd u\u017eivatel m\u016f\u017ee nav\u00edc schvalovat po\u017eadavky na p\u0159id\u00e1n\u00ed postav. Moder\u00e1tor m\u016f\u017ee nav\u00edc postavy upravovat a mazat. Administr\u00e1tor m\u016f\u017ee nav\u00edc spravovat u\u017eivatele (mazat a m\u011bnit jejich \u00farovn\u011b).\"\n}\n\n# 2) \"Aplikace je tvo\u0159ena v HTML a CSS, ...\" -> split into two paragraphs\n#    describing the Vue.js frontend / Java backend / Pinia auth.\n#    (The target paragraph has several runs split by proofErr markers, so we\n#    grab a stable Start/End snapshot first and re-derive a fresh Range from\n#    those offsets for each subsequent operation.)\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Aplikace je tvo\u0159ena v HTML a CSS, na frontend vyu\u017e\u00edv\u00e1 jazyk JavaScript a backend je v jazyku Java. Aplikace tak\u00e9 vyu\u017e\u00edv\u00e1 vlastn\u00ed MySQL datab\u00e1ze. Pro CSS je pou\u017eit framework Bootstrap.\"\n$found = $find.Execute()\nif ($found) {\n    $s = $find.Parent.Start\n    $e = $find.Parent.End\n    $afterRng = $d.Range($s, $e)\n    $afterRng.InsertAfter(\"`rNa spr\u00e1vu p\u0159ihl\u00e1\u0161en\u00ed je pou\u017eit Pinia s JWT tokeny.\")\n    $mainRng = $d.Range($s, $e)\n    $mainRng.Text = \"Frontend aplikace je tvo\u0159en pomoc\u00ed vue.js s pou\u017eit\u00edm frameworku bootstrap. Backend aplikace je v jazyce Java. Aplikace tak\u00e9 vyu\u017e\u00edv\u00e1 vlastn\u00ed MySQL datab\u00e1ze.\"\n}\n\n# 3) \"Administr\u00e1to\u0159i mohou schvalovat/zam\u00edtnout ...\" -> now also ov\u011b\u0159en\u00ed\n#    u\u017eivatel\u00e9 a moder\u00e1to\u0159i.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Administr\u00e1to\u0159i mohou schvalovat/zam\u00edtnout jednotliv\u00e9 po\u017eadavky na p\u0159id\u00e1n\u00ed postav od u\u017eivatel\u016f.\"\n$found = $find.Execute()\nif ($found) {\n    $rng = $find.Parent\n    $rng.Text = \"Administr\u00e1to\u0159i, ov\u011b\u0159en\u00ed u\u017eivatel\u00e9 a moder\u00e1to\u0159i mohou schvalovat/zam\u00edtnout jednotliv\u00e9 po\u017eadavky na p\u0159id\u00e1n\u00ed postav od u\u017eivatel\u016f.\"\n}\n\n# 4) \"Administr\u00e1to\u0159i mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je\n#    mazat.\" -> now also moder\u00e1to\u0159i, plus a brand-new following paragraph\n#    about administrators managing users.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Administr\u00e1to\u0159i mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je mazat.\"\n$found = $find.Execute()\nif ($found) {\n    $s = $find.Parent.Start\n    $e = $find.Parent.End\n    $afterRng = $d.Range($s, $e)\n    $afterRng.InsertAfter(\"`rAdministr\u00e1to\u0159i mohou upravit a mazat jednotliv\u00e9 u\u017eivatele.\")\n    $mainRng = $d.Range($s, $e)\n    $mainRng.Text = \"Administr\u00e1to\u0159i a moder\u00e1to\u0159i mohou upravovat informace o postav\u00e1ch a p\u0159\u00edpadn\u011b je mazat.\"\n}\n"}
